$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.687.60"
$ws.Range("E2").Value = "  -3.10%  "
$ws.Range("D3").Value = "1.764.94"
$ws.Range("E3").Value = "  -4.12%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.66%  "
$ws.Range("D5").Value = "241.42"
$ws.Range("E5").Value = "  -7.37%  "
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").Value = "0.4849"
$ws.Range("E7").Value = "  -7.87%  "
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "42.36"
$ws.Range("E8").Value = "  -5.51%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.2458"
$ws.Range("E9").Value = "  -22.83%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.06126"
$ws.Range("E10").Value = "  -9.79%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.787.52"
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.06559"
$ws.Range("E12").Value = "  -15.42%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "14.42"
$ws.Range("E13").Value = "  -23.04%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.6015"
$ws.Range("E14").Value = "  -23.29%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "78.14"
$ws.Range("E15").Value = "  -11.12%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "4.339"
$ws.Range("E16").Value = "  -13.48%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "25.720.60"
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "11.15"
$ws.Range("E20").Value = "  -19.43%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.026.20"
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "0.000006227"
$ws.Range("E22").Value = "  -21.72%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "3.879"
$ws.Range("E23").Value = "  -16.04%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "5.138"
$ws.Range("E24").Value = "  -14.07%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "7.981"
$ws.Range("E25").Value = "  -14.33%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "130.93"
$ws.Range("E26").Value = "  -8.04%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "1.879"
$ws.Range("E27").Value = "  -15.06%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "14.43"
$ws.Range("E28").Value = "  -14.80%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "1.357"
$ws.Range("E29").Value = "  -19.16%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "98.39"
$ws.Range("E30").Value = "  -12.09%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.08188"
$ws.Range("E31").Value = "  -6.05%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "3.576"
$ws.Range("E32").Value = "  -14.20%  "
$ws.Range("B33").Value = "Frax"
$ws.Range("C33").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D33").Value = "1.005"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.727"
$ws.Range("E34").Value = "  -4.56%  "
$ws.Range("D35").Value = "3.158"
$ws.Range("E35").Value = "  -22.55%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.04264"
$ws.Range("E36").Value = "  -12.80%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "1.019"
$ws.Range("E37").Value = "  -10.12%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.6159"
$ws.Range("E38").Value = "  -15.19%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.660"
$ws.Range("E39").Value = "  -13.98%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "2.087"
$ws.Range("E40").Value = "  -6.75%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "1.007"
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "100.63"
$ws.Range("E42").Value = "  -8.32%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.01435"
$ws.Range("E43").Value = "  -18.18%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.7839"
$ws.Range("E44").Value = "  -12.46%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.3839"
$ws.Range("E45").Value = "  -20.33%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "5.154"
$ws.Range("E46").Value = "  -12.94%  "
$ws.Range("D47").Value = "6.105"
$ws.Range("E47").Value = "  -20.13%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05176"
$ws.Range("E48").Value = "  -11.45%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "52.42"
$ws.Range("E49").Value = "  -12.18%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "1.008"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.27%  "
